$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 57
$prev = $row - 1

# New match row appended at the bottom of the results table.
$ws.Cells.Item($row, 1).Value = 56
$ws.Cells.Item($row, 2).Value = "cambodia"
$ws.Cells.Item($row, 3).Value = "cpl"
$ws.Cells.Item($row, 4).Value = "2023-2024"
$ws.Cells.Item($row, 5).Value = 45262.39583333334
$ws.Cells.Item($row, 6).Value = "Kirivong Sok Sen Chey"
$ws.Cells.Item($row, 7).Value = 1
$ws.Cells.Item($row, 8).Value = "Angkor Tiger"
$ws.Cells.Item($row, 9).Value = 3
$ws.Cells.Item($row, 10).Value = 1.91
$ws.Cells.Item($row, 11).Value = "01/12/2023 22:42"
$ws.Cells.Item($row, 12).Value = 2.09
$ws.Cells.Item($row, 13).Value = "02/12/2023 09:05"
$ws.Cells.Item($row, 14).Value = 3.58
$ws.Cells.Item($row, 15).Value = "01/12/2023 22:42"
$ws.Cells.Item($row, 16).Value = 3.64
$ws.Cells.Item($row, 17).Value = "02/12/2023 08:38"
$ws.Cells.Item($row, 18).Value = 3.22
$ws.Cells.Item($row, 19).Value = "01/12/2023 22:42"
$ws.Cells.Item($row, 20).Value = 2.83
$ws.Cells.Item($row, 21).Value = "02/12/2023 09:05"
$ws.Cells.Item($row, 22).Value = "https://www.betexplorer.com/football/cambodia/cpl/kirivong-sok-sen-chey-angkor-tiger/hpEUiaxl/"

# Carry over the formatting used by the row above (column A keeps the
# bold/bordered "Indice" look, column E keeps the date number format);
# the rest of the columns use the workbook default format, same as every
# other data row.
$ws.Range("A$prev").Copy()
$ws.Range("A$row").PasteSpecial(-4122)

$ws.Range("E$prev").Copy()
$ws.Range("E$row").PasteSpecial(-4122)

$excel.CutCopyMode = $false
